# Presupuestos.xlsx edit script
# Implements the changes described by the commit:
#   "Subo lo de taller arreglado y la just de SO"
#
# Summary of data changes (see commit diff):
#  - Economico block: Cableado quantity 2 -> 3 (row 22); new "Red Hat" USD line
#    added right after the Informix line (fills the existing blank row 37).
#  - Recomendado block: Cableado quantity 2 -> 3 (row 60); new "Patchera" line
#    added at the end of the "Red" sub-section; new "Red Hat" line added at the
#    end of the "Extra" sub-section.
#  - Premium block: Cableado quantity 2 -> 3 (old row 98); new "Patchera" line
#    added at the end of the "Red" sub-section; the existing "Red Hat" line's
#    price goes from 349 to 1300.
#  - The three section subtotal formulas (Economico / Recomendado / Premium)
#    are widened to keep covering the newly inserted rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Economico block (rows 1-37): Cableado 2 -> 3, then fill in the already
#    blank row 37 with a new "Red Hat" (USD) line. Row 38 stays blank, so no
#    rows need to be shifted for this block.
# ---------------------------------------------------------------------------
$ws.Range("E22").Value = 3

$ws.Range("C37").Value = "Red Hat"
$ws.Range("D37").Value = 1300
$ws.Range("E37").Value = 1
$ws.Range("F37").Value = 1300
$ws.Range("G37").Value = "USD"

# ---------------------------------------------------------------------------
# 2) Recomendado block: Cableado 2 -> 3 (still row 60 at this point, nothing
#    shifted yet).
# ---------------------------------------------------------------------------
$ws.Range("E60").Value = 3

# Insert the new "Patchera" line at the end of the "Red" sub-section (right
# before the "Servidor" header, currently row 64). This pushes everything
# from the old row 64 onward down by one row.
$ws.Rows(63).Insert()
$ws.Range("C63").Value = "Patchera"
$ws.Range("D63").Value = 90
$ws.Range("E63").Value = 8
$ws.Range("F63").Formula = "=D63*E63"
$ws.Range("G63").Value = "USD"

# Insert the new "Red Hat" line at the end of the "Extra" sub-section (the
# blank separator row that currently sits right before the "Premium" header,
# now at row 76 after the previous insert). This pushes the "Premium" header
# and everything below it down by one more row.
$ws.Rows(76).Insert()
$ws.Range("C76").Value = "Red Hat"
$ws.Range("D76").Value = 1300
$ws.Range("E76").Value = 1
$ws.Range("F76").Value = 1300
$ws.Range("G76").Value = "USD"

# ---------------------------------------------------------------------------
# 3) Premium block: Cableado 2 -> 3. After the two inserts above, the old
#    row 98 now sits at row 100.
# ---------------------------------------------------------------------------
$ws.Range("E100").Value = 3

# Insert the new "Patchera" line at the end of the "Red" sub-section of the
# Premium block (right before the "Servidor" header, now at row 102 after
# the earlier shifts).
$ws.Rows(103).Insert()
$ws.Range("C103").Value = "Patchera"
$ws.Range("D103").Value = 90
$ws.Range("E103").Value = 8
$ws.Range("F103").Formula = "=D103*E103"
$ws.Range("G103").Value = "USD"

# Existing "Red Hat" line in the Premium/Extra sub-section: price 349 -> 1300.
# After all the inserts above, the old row 113 now sits at row 116.
$ws.Range("D116").Value = 1300

# ---------------------------------------------------------------------------
# 4) Fix up the three section subtotal formulas so they explicitly cover the
#    newly inserted rows (set explicitly rather than relying on Excel's
#    auto-extend-on-insert heuristics, which don't apply to the boundary
#    insert positions used above).
# ---------------------------------------------------------------------------
$ws.Range("B3").Formula = "=SUM(F6:F37)"
$ws.Range("B39").Formula = "=SUM(F42:F76)"
$ws.Range("B78").Formula = "=SUM(F81:F116)"

# ---------------------------------------------------------------------------
# 5) Misc view-state bookkeeping matching the commit (dimension grows to
#    L116 automatically from the cells above; refresh the active selection).
# ---------------------------------------------------------------------------
$ws.Range("E101").Select()
